$wb = $excel.ActiveWorkbook

# --- Rename the first sheet from "贷款明细表" to "贷款信息表" ---
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "贷款信息表"

# --- Update the selection on the "LPR" sheet (sheet 3). It is currently the ---
# --- tab-selected sheet, so touch its selection first, before handing the  ---
# --- active/selected tab over to sheet 1.                                  ---
$ws3.Activate()
$ws3.Range("D6").Select() | Out-Null

# --- Make sheet 1 ("贷款信息表") the active / tab-selected sheet, and move ---
# --- its selection to E26:F26.                                             ---
$ws1.Activate()
$ws1.Range("E26:F26").Select() | Out-Null
